$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(51,4).NumberFormat = "@"

$ws.Cells.Item(2,4).Value = "25.900.57"
$ws.Cells.Item(2,5).Value = "  +0.33%  "
$ws.Cells.Item(3,4).Value = "1.647.55"
$ws.Cells.Item(3,5).Value = "  +0.76%  "
$ws.Cells.Item(4,4).Value = "1.007"
$ws.Cells.Item(4,5).Value = "  +0.49%  "
$ws.Cells.Item(5,4).Value = "215.38"
$ws.Cells.Item(5,5).Value = "  -0.04%  "
$ws.Cells.Item(6,4).Value = "0.5079"
$ws.Cells.Item(6,5).Value = "  +1.00%  "
$ws.Cells.Item(7,4).Value = "1.006"
$ws.Cells.Item(7,5).Value = "  +0.37%  "
$ws.Cells.Item(8,4).Value = "0.2574"
$ws.Cells.Item(8,5).Value = "  -0.10%  "
$ws.Cells.Item(9,4).Value = "0.06418"
$ws.Cells.Item(9,5).Value = "  +0.05%  "
$ws.Cells.Item(10,4).Value = "19.72"
$ws.Cells.Item(10,5).Value = "  -0.08%  "
$ws.Cells.Item(11,4).Value = "0.07774"
$ws.Cells.Item(11,5).Value = "  +1.01%  "
$ws.Cells.Item(12,4).Value = "4.310"
$ws.Cells.Item(12,5).Value = "  +1.24%  "
$ws.Cells.Item(13,4).Value = "1.639.43"
$ws.Cells.Item(13,5).Value = "  +0.25%  "
$ws.Cells.Item(14,4).Value = "0.5469"
$ws.Cells.Item(14,5).Value = "  +0.16%  "
$ws.Cells.Item(15,4).Value = "0.0₅7905"
$ws.Cells.Item(15,5).Value = "  -0.36%  "
$ws.Cells.Item(16,4).Value = "65.18"
$ws.Cells.Item(16,5).Value = "  +2.48%  "
$ws.Cells.Item(17,4).Value = "26.000.50"
$ws.Cells.Item(17,5).Value = "  +0.62%  "
$ws.Cells.Item(18,4).Value = "1.006"
$ws.Cells.Item(18,5).Value = "  +0.39%  "
$ws.Cells.Item(19,4).Value = "197.50"
$ws.Cells.Item(19,5).Value = "  -2.91%  "
$ws.Cells.Item(20,4).Value = "4.422"
$ws.Cells.Item(20,5).Value = "  +2.15%  "
$ws.Cells.Item(21,4).Value = "10.03"
$ws.Cells.Item(21,5).Value = "  +0.82%  "
$ws.Cells.Item(22,4).Value = "6.075"
$ws.Cells.Item(22,5).Value = "  +1.68%  "
$ws.Cells.Item(23,4).Value = "1.009"
$ws.Cells.Item(23,5).Value = "  +0.53%  "
$ws.Cells.Item(24,4).Value = "1.862"
$ws.Cells.Item(25,4).Value = "141.27"
$ws.Cells.Item(25,5).Value = "  -0.05%  "
$ws.Cells.Item(26,4).Value = "0.1145"
$ws.Cells.Item(26,5).Value = "  -0.11%  "
$ws.Cells.Item(27,4).Value = "6.906"
$ws.Cells.Item(27,5).Value = "  +2.90%  "
$ws.Cells.Item(28,4).Value = "15.74"
$ws.Cells.Item(28,5).Value = "  +0.23%  "
$ws.Cells.Item(29,4).Value = "1.243"
$ws.Cells.Item(29,5).Value = "  +0.10%  "
$ws.Cells.Item(30,4).Value = "0.05029"
$ws.Cells.Item(30,5).Value = "  +0.15%  "
$ws.Cells.Item(31,4).Value = "3.277"
$ws.Cells.Item(31,5).Value = "  +0.13%  "
$ws.Cells.Item(32,4).Value = "3.206"
$ws.Cells.Item(32,5).Value = "  +0.56%  "
$ws.Cells.Item(33,4).Value = "1.541"
$ws.Cells.Item(33,5).Value = "  +0.12%  "
$ws.Cells.Item(34,4).Value = "2.369"
$ws.Cells.Item(34,5).Value = "  +0.69%  "
$ws.Cells.Item(35,4).Value = "0.8934"
$ws.Cells.Item(35,5).Value = "  -0.03%  "
$ws.Cells.Item(36,4).Value = "2.595"
$ws.Cells.Item(36,5).Value = "  -0.53%  "
$ws.Cells.Item(37,4).Value = "1.134.03"
$ws.Cells.Item(37,5).Value = "  -3.20%  "
$ws.Cells.Item(38,4).Value = "0.5539"
$ws.Cells.Item(38,5).Value = "  -1.02%  "
$ws.Cells.Item(39,4).Value = "0.01564"
$ws.Cells.Item(39,5).Value = "  +0.20%  "
$ws.Cells.Item(40,4).Value = "1.008"
$ws.Cells.Item(40,5).Value = "  +0.53%  "
$ws.Cells.Item(41,4).Value = "5.670"
$ws.Cells.Item(41,5).Value = "  -0.03%  "
$ws.Cells.Item(42,4).Value = "0.8154"
$ws.Cells.Item(42,5).Value = "  +0.86%  "
$ws.Cells.Item(43,4).Value = "99.81"
$ws.Cells.Item(43,5).Value = "  +0.32%  "
$ws.Cells.Item(44,4).Value = "0.0₈122"
$ws.Cells.Item(44,5).Value = "  +6.31%  "
$ws.Cells.Item(45,4).Value = "1.784.98"
$ws.Cells.Item(45,5).Value = "  +0.76%  "
$ws.Cells.Item(46,4).Value = "0.4540"
$ws.Cells.Item(46,5).Value = "  +0.58%  "
$ws.Cells.Item(47,2).Value = "Frax"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(47,4).Value = "1.007"
$ws.Cells.Item(47,5).Value = "  +0.34%  "
$ws.Cells.Item(48,2).Value = "Aave"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48,4).Value = "55.21"
$ws.Cells.Item(48,5).Value = "  +0.48%  "
$ws.Cells.Item(49,4).Value = "0.05094"
$ws.Cells.Item(49,5).Value = "  +1.09%  "
$ws.Cells.Item(50,2).Value = "USDD"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(50,4).Value = "1.009"
$ws.Cells.Item(50,5).Value = "  +0.43%  "
$ws.Cells.Item(51,2).Value = "Algorand"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51,4).Value = "0.09563"
$ws.Cells.Item(51,5).Value = "  +3.02%  "
